$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.474.43'
$ws.Range("E2").Value = '  -2.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.743.08'
$ws.Range("E3").Value = '  -3.36%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.72'
$ws.Range("E5").Value = '  -4.34%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4209'
$ws.Range("E7").Value = '  -9.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3570'
$ws.Range("E8").Value = '  -5.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.58'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07408'
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.111'
$ws.Range("E11").Value = '  -3.41%  '
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.43'
$ws.Range("E13").Value = '  -4.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.098'
$ws.Range("E14").Value = '  -3.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.170'
$ws.Range("E15").Value = '  -3.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.740.75'
$ws.Range("E16").Value = '  -3.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001063'
$ws.Range("E17").Value = '  -2.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.46'
$ws.Range("E18").Value = '  +6.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06164'
$ws.Range("E19").Value = '  -8.57%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.86'
$ws.Range("E21").Value = '  -3.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.100'
$ws.Range("E22").Value = '  -4.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5249'
$ws.Range("E23").Value = '  -5.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.498.60'
$ws.Range("E24").Value = '  -2.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.58'
$ws.Range("E25").Value = '  -2.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.322'
$ws.Range("E26").Value = '  -3.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.42'
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.37'
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.363'
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.938.65'
$ws.Range("E30").Value = '  -3.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '125.91'
$ws.Range("E31").Value = '  -5.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.203'
$ws.Range("E32").Value = '  -4.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.684'
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09140'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.687'
$ws.Range("E35").Value = '  -8.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.64'
$ws.Range("E36").Value = '  +4.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02286'
$ws.Range("E37").Value = '  -2.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2125'
$ws.Range("E38").Value = '  -4.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.076'
$ws.Range("E39").Value = '  -3.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06071'
$ws.Range("E40").Value = '  -4.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6375'
$ws.Range("E41").Value = '  -3.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.188'
$ws.Range("E42").Value = '  -3.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.423'
$ws.Range("E43").Value = '  -4.96%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.889'
$ws.Range("E45").Value = '  -4.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.75'
$ws.Range("E46").Value = '  -3.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.715'
$ws.Range("E47").Value = '  -3.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5850'
$ws.Range("E48").Value = '  -4.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.01'
$ws.Range("E49").Value = '  -3.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.948'
$ws.Range("E50").Value = '  -4.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06847'
$ws.Range("E51").Value = '  -4.45%  '
